$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows at the top of the data table (above old row 2), shifting
# all existing data + formulas down. Excel auto-adjusts the relative row
# formulas in column A as part of the insert.
$ws.Rows("2:8").Insert()

# The new rows (2-8) need their date formulas filled in, continuing the
# existing "(row below) + 1" pattern that row 9 (formerly row 2) already
# carries after the shift.
for ($r = 8; $r -ge 2; $r--) {
    $below = $r + 1
    $ws.Range("A$r").Formula = "=(A$below+1)"
}

# Full refreshed dataset (dates already correct via formulas above / the
# shifted existing formulas): set Lower CI / Upper CI / Estimate of Rt for
# every data row 2-58.
$rowsData = @(
    @(2, 0.62, 1.08, 0.85),
    @(3, 0.64, 1.07, 0.85),
    @(4, 0.66, 1.05, 0.85),
    @(5, 0.67, 1.04, 0.86),
    @(6, 0.69, 1.03, 0.87),
    @(7, 0.71, 1.02, 0.87),
    @(8, 0.73, 1.02, 0.88),
    @(9, 0.75, 1.01, 0.89),
    @(10, 0.77, 1.01, 0.9),
    @(11, 0.8, 1, 0.91),
    @(12, 0.82, 1, 0.92),
    @(13, 0.84, 1.01, 0.93),
    @(14, 0.87, 1.01, 0.94),
    @(15, 0.89, 1.01, 0.95),
    @(16, 0.91, 1.02, 0.97),
    @(17, 0.93, 1.04, 0.98),
    @(18, 0.94, 1.05, 1),
    @(19, 0.96, 1.06, 1.01),
    @(20, 0.98, 1.08, 1.03),
    @(21, 0.99, 1.1000000000000001, 1.04),
    @(22, 1.01, 1.1200000000000001, 1.05),
    @(23, 1.02, 1.1299999999999999, 1.07),
    @(24, 1.04, 1.1499999999999999, 1.08),
    @(25, 1.05, 1.1599999999999999, 1.1000000000000001),
    @(26, 1.06, 1.17, 1.1100000000000001),
    @(27, 1.08, 1.18, 1.1299999999999999),
    @(28, 1.0900000000000001, 1.18, 1.1299999999999999),
    @(29, 1.1000000000000001, 1.19, 1.1399999999999999),
    @(30, 1.1000000000000001, 1.2, 1.1499999999999999),
    @(31, 1.1100000000000001, 1.21, 1.1599999999999999),
    @(32, 1.1200000000000001, 1.23, 1.17),
    @(33, 1.1299999999999999, 1.23, 1.18),
    @(34, 1.1299999999999999, 1.24, 1.19),
    @(35, 1.1399999999999999, 1.25, 1.2),
    @(36, 1.1399999999999999, 1.26, 1.2),
    @(37, 1.1599999999999999, 1.27, 1.21),
    @(38, 1.17, 1.27, 1.22),
    @(39, 1.18, 1.28, 1.23),
    @(40, 1.18, 1.29, 1.23),
    @(41, 1.19, 1.3, 1.24),
    @(42, 1.19, 1.31, 1.24),
    @(43, 1.19, 1.32, 1.24),
    @(44, 1.19, 1.32, 1.24),
    @(45, 1.19, 1.32, 1.24),
    @(46, 1.19, 1.32, 1.24),
    @(47, 1.18, 1.31, 1.23),
    @(48, 1.18, 1.3, 1.23),
    @(49, 1.17, 1.29, 1.22),
    @(50, 1.1599999999999999, 1.28, 1.21),
    @(51, 1.1599999999999999, 1.26, 1.2),
    @(52, 1.1499999999999999, 1.25, 1.19),
    @(53, 1.1299999999999999, 1.23, 1.18),
    @(54, 1.1200000000000001, 1.22, 1.17),
    @(55, 1.1000000000000001, 1.2, 1.1599999999999999),
    @(56, 1.0900000000000001, 1.19, 1.1399999999999999),
    @(57, 1.07, 1.18, 1.1299999999999999),
    @(58, 1.06, 1.1599999999999999, 1.1200000000000001)
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
}

# Restore the author's final selection.
$ws.Range("B3").Select()
